# Weekly update: insert this week's "Cilantro" price rows at the top of the
# data block (row 189), pushing all the existing history rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows - everything currently on/after row 189
# (old rows 189:236) shifts down to 191:238.
$ws.Rows("189:190").Insert()

# Row 189 - "Primera" quality
$ws.Range("A189").Value = 11
$ws.Range("B189").Value = 'Vega Monumental Concepción'
$ws.Range("C189").Value = 'Bíobío'
$ws.Range("D189").Value = 44855
$ws.Range("E189").Value = 8
$ws.Range("F189").Value = 100112040
$ws.Range("G189").Value = 'Cilantro'
$ws.Range("H189").Value = 'Sin especificar'
$ws.Range("I189").Value = 'Primera'
$ws.Range("J189").Value = 200
$ws.Range("K189").Value = 700
$ws.Range("L189").Value = 800
$ws.Range("M189").Value = 750
$ws.Range("N189").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O189").Value = 'Región de Ñuble'
$ws.Range("P189").Value = 750
$ws.Range("Q189").Value = 1
$ws.Range("R189").Value = 'Hortaliza'

# Row 190 - "Segunda" quality
$ws.Range("A190").Value = 11
$ws.Range("B190").Value = 'Vega Monumental Concepción'
$ws.Range("C190").Value = 'Bíobío'
$ws.Range("D190").Value = 44855
$ws.Range("E190").Value = 8
$ws.Range("F190").Value = 100112040
$ws.Range("G190").Value = 'Cilantro'
$ws.Range("H190").Value = 'Sin especificar'
$ws.Range("I190").Value = 'Segunda'
$ws.Range("J190").Value = 100
$ws.Range("K190").Value = 600
$ws.Range("L190").Value = 600
$ws.Range("M190").Value = 600
$ws.Range("N190").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O190").Value = 'Región de Ñuble'
$ws.Range("P190").Value = 600
$ws.Range("Q190").Value = 1
$ws.Range("R190").Value = 'Hortaliza'
